# Fruta / hortaliza, semanal
# Insert a new record row at row 23 (pushing the existing rows 23..87 down
# to 24..88), then populate the new row's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 23..87 down by one to make room for the new record.
$ws.Rows("23:23").Insert()

# Populate the new row 23 with the inserted record's values. The columns
# that are constant across the whole sheet (A,B,C,E,F,G,H,I,J,K,Q,R,T) are
# copied from the (now shifted) neighbouring row 24; the rest get the new
# record's own values.
$ws.Range("A23").Value = 10
$ws.Range("B23").Value = "Vega Modelo de Temuco"
$ws.Range("C23").Value = "La Araucanía"
$ws.Range("D23").Value = 44498
$ws.Range("D23").NumberFormat = $ws.Range("D24").NumberFormat
$ws.Range("E23").Value = 9
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100107
$ws.Range("H23").Value = "Otros"
$ws.Range("I23").Value = 100107002
$ws.Range("J23").Value = "Chirimoya"
$ws.Range("K23").Value = "Cultivar IV Región"
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 300
$ws.Range("N23").Value = 3000
$ws.Range("O23").Value = 3000
$ws.Range("P23").Value = 3000
$ws.Range("Q23").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R23").Value = "Provincia del Elquí"
$ws.Range("S23").Value = 3000
$ws.Range("T23").Value = 1
